$d = $word.ActiveDocument

# 1. Refresh the generation timestamp stamped in the footer.
$footer = $d.Sections(1).Footers(1).Range
$footer.Find.Execute("2025-06-30 12:13Z / ", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "2025-07-02 02:48Z / ", 2)

# 2. Add the regression-test character styles (b / i / sub / sup / u),
#    mirroring the PubMed inline-markup shorthand styles.
$wdStyleTypeCharacter = 2
$wdUnderlineSingle = 1

$b = $d.Styles.Add("b", $wdStyleTypeCharacter)
$b.BaseStyle = "DefaultParagraphFont"
$b.Priority = 1
$b.QuickStyle = $true
$b.Font.Bold = $true

$i = $d.Styles.Add("i", $wdStyleTypeCharacter)
$i.BaseStyle = "DefaultParagraphFont"
$i.Priority = 1
$i.QuickStyle = $true
$i.Font.Italic = $true

$sub = $d.Styles.Add("sub", $wdStyleTypeCharacter)
$sub.BaseStyle = "DefaultParagraphFont"
$sub.Priority = 1
$sub.QuickStyle = $true
$sub.Font.Subscript = $true

$sup = $d.Styles.Add("sup", $wdStyleTypeCharacter)
$sup.BaseStyle = "DefaultParagraphFont"
$sup.Priority = 1
$sup.QuickStyle = $true
$sup.Font.Superscript = $true

$u = $d.Styles.Add("u", $wdStyleTypeCharacter)
$u.BaseStyle = "DefaultParagraphFont"
$u.Priority = 1
$u.QuickStyle = $true
$u.Font.Underline = $wdUnderlineSingle
